$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.446.15"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "2.272.79"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "502.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.527"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.39%  "
$ws.Range("D9").Value = "2.285.11"
$ws.Range("E9").Value = "  -0.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0961"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.55%  "
$ws.Range("E11").Value = "  +1.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.340"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.90"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.87%  "
$ws.Range("D15").Value = "2.677.93"
$ws.Range("E15").Value = "  -0.27%  "
$ws.Range("D16").Value = "54.534.00"
$ws.Range("E16").Value = "  +0.46%  "
$ws.Range("E17").Value = "  +0.67%  "
$ws.Range("D18").Value = "2.286.62"
$ws.Range("E18").Value = "  -1.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.30"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "305.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.42"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.995"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.33%  "
$ws.Range("E26").Value = "  -1.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.42"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "170.53"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.04"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.69%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.62"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.29%  "
$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").Value = "0.0₃0697"
$ws.Range("E31").Value = "  +1.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.11"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.12%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.89"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.996"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.906"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.91%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.76"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.42"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.80%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.373"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.03"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.38"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "125.94"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0494"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.89%  "
$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "248.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0900"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.545"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.13%  "
$ws.Range("E49").Value = "  +0.18%  "
$ws.Range("E50").Value = "  +0.46%  "
$ws.Range("E51").Value = "  +0.43%  "
